$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDNS")

# Row 4: Inventory
$ws.Range("B4").Value = 77000000.0
$ws.Range("C4").Value = 76000000.0
$ws.Range("D4").Value = 48000000.0
$ws.Range("E4").Value = 44000000.0
$ws.Range("F4").Value = 61000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 268000000.0
$ws.Range("C15").Value = 341000000.0
$ws.Range("D15").Value = 290000000.0
$ws.Range("E15").Value = 306000000.0
$ws.Range("F15").Value = 233000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -729000000.0
$ws.Range("C22").Value = -732000000.0
$ws.Range("D22").Value = -723000000.0
$ws.Range("E22").Value = -705000000.0
$ws.Range("F22").Value = -699000000.0

# Row 34: Net Debt
$ws.Range("G34").Value = -333191000.0

# Row 35: Total Debt
$ws.Range("G35").Value = 372019000.0
